# Updates cryptos list values (price/volume columns) to the latest scrape.
# Generated from the authoritative cell-level diff between before/after OOXML.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.830.15"
$ws.Range("E2").Value = "  -1.27%  "

# Row 3
$ws.Range("D3").Value = "1.633.39"
$ws.Range("E3").Value = "  -1.29%  "

# Row 4
$ws.Range("E4").Value = "  -0.38%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.81"
$ws.Range("E5").Value = "  -0.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5015"
$ws.Range("E6").Value = "  -1.98%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2563"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06393"
$ws.Range("E9").Value = "  -0.14%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.58"
$ws.Range("E10").Value = "  -1.76%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07672"
$ws.Range("E11").Value = "  -1.80%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.634.88"
$ws.Range("E12").Value = "  -1.28%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.238"
$ws.Range("E13").Value = "  -0.97%  "

# Row 14
$ws.Range("D14").Value = "1.858.56"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5417"
$ws.Range("E15").Value = "  -1.78%  "

# Row 16
$ws.Range("E16").Value = "  -1.22%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.45"
$ws.Range("E17").Value = "  -0.61%  "

# Row 18
$ws.Range("D18").Value = "25.843.24"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.39%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.88"
$ws.Range("E20").Value = "  -3.18%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.323"
$ws.Range("E21").Value = "  -2.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.914"
$ws.Range("E22").Value = "  -1.45%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.957"
$ws.Range("E23").Value = "  -0.94%  "

# Row 24
$ws.Range("E24").Value = "  -0.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.907"
$ws.Range("E25").Value = "  +10.61%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.24"
$ws.Range("E26").Value = "  -1.37%  "

# Row 27
$ws.Range("E27").Value = "  -2.38%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.65"
$ws.Range("E28").Value = "  -0.76%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.695"
$ws.Range("E29").Value = "  -3.93%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.239"
$ws.Range("E30").Value = "  -0.17%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04982"
$ws.Range("E31").Value = "  -2.86%  "

# Row 32
$ws.Range("E32").Value = "  -2.56%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.172"
$ws.Range("E33").Value = "  -1.38%  "

# Row 34
$ws.Range("E34").Value = "  -1.62%  "

# Row 35
$ws.Range("E35").Value = "  -0.41%  "

# Row 36
$ws.Range("D36").Value = "1.169.55"
$ws.Range("E36").Value = "  +1.19%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8901"
$ws.Range("E37").Value = "  -4.08%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.612"
$ws.Range("E38").Value = "  -4.96%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5577"
$ws.Range("E39").Value = "  -1.88%  "

# Row 40
$ws.Range("E40").Value = "  -2.04%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.555"
$ws.Range("E41").Value = "  -0.27%  "

# Row 42
$ws.Range("E42").Value = "  -0.37%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.679"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8071"
$ws.Range("E44").Value = "  -3.20%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.31"
$ws.Range("E45").Value = "  -0.84%  "

# Row 46
$ws.Range("D46").Value = "1.770.32"
$ws.Range("E46").Value = "  -1.20%  "

# Row 47
$ws.Range("E47").Value = "  -1.62%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4514"
$ws.Range("E48").Value = "  -0.72%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  -0.33%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.59"
$ws.Range("E50").Value = "  -1.93%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05075"
$ws.Range("E51").Value = "  +0.77%  "
